# Append a new results row (row 7) to Sheet1, extending the data range
# from A1:J6 to A1:J7, matching a new run with B=50 (n_points) for A=2 (degree).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 50
$ws.Range("C7").Value = 0.00623347282409668
$ws.Range("D7").Value = 1200.043860697746
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "-"
$ws.Range("G7").Value = 306.27301
$ws.Range("H7").Value = 85.36555
$ws.Range("I7").Value = "-"
$ws.Range("J7").Value = "-"
